$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.282.76"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "3.497.10"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'586.18"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "'134.24"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("E11").Value = "  +2.54%  "
$ws.Range("D12").Value = "4.094.91"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "3.499.90"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'25.83"
$ws.Range("E16").Value = "  -5.63%  "
$ws.Range("D17").Value = "64.286.20"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "'9.87"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("D21").Value = "'393.54"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "3.638.09"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'74.36"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "'1.50"
$ws.Range("E29").Value = "  -3.68%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.24"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.27"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "3.518.09"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  +3.93%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "'23.42"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "'5.15"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").Value = "'165.87"
$ws.Range("E39").Value = "  +4.56%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "'25.23"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").Value = "2.461.67"
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("D48").Value = "'6.77"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "'21.14"
$ws.Range("E51").Value = "  -0.41%  "

# Restore default (unstyled) appearance for text-forced numeric-looking cells
foreach ($addr in @("D5","D6","D16","D18","D21","D24","D27","D29","D30","D31","D35","D36","D39","D43","D48","D51")) {
    $ws.Range($addr).Style = "Normal"
}
